$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: price input 29 -> 33 ---
$ws1.Range("E2").Value = 33

# --- Sheet2, row 23: decay multiplier 0.93 -> 0.95 ---
$ws2.Range("K23").Formula = "=J23*0.95"
$ws2.Range("L23").Formula = "=K23*0.95"
$ws2.Range("M23").Formula = "=L23*0.95"
$ws2.Range("N23").Formula = "=M23*0.95"
$ws2.Range("O23").Formula = "=N23*0.95"

# --- Sheet2, row 26: newly populated ratio row ---
$ws2.Range("H26").Formula = "=H27/H2"
$ws2.Range("I26").Formula = "=I27/I2"
$ws2.Range("J26").Formula = "=J27/J2"
$ws2.Range("K26").Formula = "=J26*1.05"
$ws2.Range("L26").Formula = "=K26*1.05"
$ws2.Range("M26").Formula = "=L26*1.05"
$ws2.Range("N26").Formula = "=M26*1.05"
$ws2.Range("O26").Formula = "=N26*1.05"

# --- Sheet2, row 27: newly populated forecast row ---
$ws2.Range("K27").Formula = "=K26*K2"
$ws2.Range("L27").Formula = "=L26*L2"
$ws2.Range("M27").Formula = "=M26*M2"
$ws2.Range("N27").Formula = "=N26*N2"
$ws2.Range("O27").Formula = "=O26*O2"

# --- Sheet2, row 28: newly populated forecast row ---
$ws2.Range("K28").Formula = "=K2*0.02"
$ws2.Range("L28").Formula = "=L2*0.02"
$ws2.Range("M28").Formula = "=M2*0.02"
$ws2.Range("N28").Formula = "=N2*0.02"
$ws2.Range("O28").Formula = "=O2*0.02"

# --- Sheet2, row 29: now derived from rows 27/28 rather than row 22 ---
$ws2.Range("K29").Formula = "=K27-K28"
$ws2.Range("L29").Formula = "=L27-L28"
$ws2.Range("M29").Formula = "=M27-M28"
$ws2.Range("N29").Formula = "=N27-N28"
$ws2.Range("O29").Formula = "=O27-O28"

# --- Sheet2, row 35: NPV now based on row 29 instead of row 14 ---
$ws2.Range("R35").Formula = "=NPV(R34,K29:EQ29)+Sheet1!E5-Sheet1!E6"

# --- View state tweaks ---
$ws1.Range("E3").Select()
$ws2.Range("L23").Select()

$excel.ActiveWindow.SplitColumn = 1
$excel.ActiveWindow.SplitRow = 1
